$wb = $excel.ActiveWorkbook

# Grab reference to the existing "ODI Batting" sheet before inserting a new one.
$odiSheet = $wb.Worksheets.Item(1)

# Insert a brand-new worksheet before the ODI Batting sheet for the player info table.
$playerInfoSheet = $wb.Worksheets.Add($odiSheet)
$playerInfoSheet.Name = "Player Info"

# Re-fetch the ODI Batting sheet by name since the previous reference's position shifted.
$odiSheet = $wb.Worksheets.Item("ODI Batting")

# --- Populate the new "Player Info" sheet ---
# Copy the existing bold/bordered/centered header style from ODI Batting's header row
# onto the new header row so it re-uses the same style instead of creating a new one.
$odiSheet.Range("A1").Copy($playerInfoSheet.Range("A1:D1"))

$playerInfoSheet.Range("A1").Value = "ID"
$playerInfoSheet.Range("B1").Value = "NAME"
$playerInfoSheet.Range("C1").Value = "BATTING_HAND"
$playerInfoSheet.Range("D1").Value = "BOWL_STYLE"

$playerInfoSheet.Range("A2").NumberFormat = "@"
$playerInfoSheet.Range("A2").Value = "5860"
$playerInfoSheet.Range("B2").Value = "Joshua Ryan Philippe"
$playerInfoSheet.Range("C2").Value = "Right Handed"
$playerInfoSheet.Range("D2").Value = "Does Not Bowl | Unknown"

# --- Update the "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and
#     shrink the URLs down to just the numeric match code (kept as text) ---
$odiSheet.Range("D1").Value = "MATCH_CODE"
$odiSheet.Range("D2:D4").NumberFormat = "@"
$odiSheet.Range("D2").Value = "4483"
$odiSheet.Range("D3").Value = "4484"
$odiSheet.Range("D4").Value = "4486"
